# "Generate Report for Archive"
#
# The localization status report is regenerated: every "Ready for handoff"
# status cell becomes "In Translation", and the Status columns (which were
# sized to fit the old, longer text) are re-fitted to the new, shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The fitted width Excel computes for the new text in these Status columns
# (Calibri 11, same as the rest of the workbook).
$newStatusColWidth = 12.576851254417766

foreach ($ws in $wb.Worksheets) {

    $used = $ws.UsedRange
    $headerRow = $used.Rows.Item(1)

    # Find every "Status" header cell (col C on the per-language sheets) and
    # every per-language status header (zh-cn / de-de, cols E/F on Overview),
    # then update the matching data cells + refit that column's width.
    foreach ($headerCell in $headerRow.Cells) {
        $header = $headerCell.Value()
        if ($header -eq "Status" -or $header -eq "zh-cn" -or $header -eq "de-de") {
            $col = $headerCell.Column

            foreach ($row in $used.Rows) {
                $cell = $ws.Cells.Item($row.Row, $col)
                if ($cell.Value() -eq $oldStatus) {
                    $cell.Value = $newStatus
                }
            }

            $ws.Columns.Item($col).ColumnWidth = $newStatusColWidth
        }
    }
}
